$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 2499.125
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2499.125
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 7497.375
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -8037.375
# Row 73
$ws.Range("H73").Value = 2499.125
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2499.125
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 7497.375
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -9369.375
# Row 98
$ws.Range("H98").Value = 1691.5714
$ws.Range("I98").Value = 1691.5714
$ws.Range("K98").Value = 1691.5714
$ws.Range("M98").Value = -193.5714
# Row 99
$ws.Range("H99").Value = 3852
$ws.Range("I99").Value = 1316
$ws.Range("J99").Value = 13996
$ws.Range("K99").Value = 3948
$ws.Range("L99").Value = 41988
$ws.Range("M99").Value = -2450
$ws.Range("N99").Value = -44984
# Row 122
$ws.Range("H122").Value = 1691.5714
$ws.Range("I122").Value = 1691.5714
$ws.Range("K122").Value = 5074.7142
$ws.Range("M122").Value = -2624.7142
# Row 132
$ws.Range("H132").Value = 111263.66
$ws.Range("I132").Value = 235830.03
$ws.Range("K132").Value = 707490.09
$ws.Range("M132").Value = -704960.09
# Row 135
$ws.Range("H135").Value = 2367
$ws.Range("I135").Value = 1607
$ws.Range("K135").Value = 14463
$ws.Range("M135").Value = -11928
# Row 137
$ws.Range("H137").Value = 2803.0833
$ws.Range("I137").Value = 1873.7
$ws.Range("J137").Value = 7450
$ws.Range("K137").Value = 5621.1
$ws.Range("L137").Value = 22350
$ws.Range("M137").Value = -3071.1
$ws.Range("N137").Value = -27450
# Row 138
$ws.Range("H138").Value = 6587.6416
$ws.Range("I138").Value = 2947.6875
$ws.Range("J138").Value = 8161.676
$ws.Range("K138").Value = 8843.0625
$ws.Range("L138").Value = 24485.028
$ws.Range("M138").Value = -3703.0625
$ws.Range("N138").Value = -34765.02800000001
# Row 139
$ws.Range("H139").Value = 239999
$ws.Range("J139").Value = 239999
$ws.Range("L139").Value = 239999
$ws.Range("N139").Value = -250279
# Row 140
$ws.Range("H140").Value = 60961.11
$ws.Range("J140").Value = 59831.25
$ws.Range("L140").Value = 59831.25
$ws.Range("N140").Value = -70191.25

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3033.9443
$ws.Range("I45").Value = 1899.9
$ws.Range("J45").Value = 4451.5
$ws.Range("K45").Value = 1899.9
$ws.Range("L45").Value = 4451.5
$ws.Range("M45").Value = -1522.9
$ws.Range("N45").Value = -5205.5
# Row 61
$ws.Range("H61").Value = 10786.237
$ws.Range("I61").Value = 11040.363
$ws.Range("J61").Value = 10436.8125
$ws.Range("K61").Value = 11040.363
$ws.Range("L61").Value = 10436.8125
$ws.Range("M61").Value = -10828.363
$ws.Range("N61").Value = -10860.8125
# Row 122
$ws.Range("H122").Value = 310205.22
$ws.Range("I122").Value = 502315.38
$ws.Range("K122").Value = 1506946.14
$ws.Range("M122").Value = -1504496.14
# Row 131
$ws.Range("H131").Value = 92975
$ws.Range("J131").Value = 92975
$ws.Range("L131").Value = 92975
$ws.Range("N131").Value = -103055
# Row 132
$ws.Range("H132").Value = 12467.396
$ws.Range("I132").Value = 16060
$ws.Range("J132").Value = 5772.091
$ws.Range("K132").Value = 48180
$ws.Range("L132").Value = 17316.273
$ws.Range("M132").Value = -45650
$ws.Range("N132").Value = -22376.273
# Row 136
$ws.Range("H136").Value = 10786.237
$ws.Range("I136").Value = 11040.363
$ws.Range("J136").Value = 10436.8125
$ws.Range("K136").Value = 33121.089
$ws.Range("L136").Value = 31310.4375
$ws.Range("M136").Value = -30571.089
$ws.Range("N136").Value = -36410.4375

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3037.4717
$ws.Range("I134").Value = 2343.122
$ws.Range("J134").Value = 5409.8335
$ws.Range("K134").Value = 7029.366
$ws.Range("L134").Value = 16229.5005
$ws.Range("M134").Value = -4494.366
$ws.Range("N134").Value = -21299.5005

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 875
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -726
# Row 58
$ws.Range("H58").Value = 1431841.8
$ws.Range("I58").Value = 2502699.2
$ws.Range("K58").Value = 2502699.2
$ws.Range("M58").Value = -2502496.2
# Row 132
$ws.Range("H132").Value = 74087780
$ws.Range("I132").Value = 83336250
$ws.Range("J132").Value = 99999.5
$ws.Range("K132").Value = 250008750
$ws.Range("L132").Value = 299998.5
$ws.Range("M132").Value = -250006220
$ws.Range("N132").Value = -305058.5
# Row 134
$ws.Range("H134").Value = 770.46155
$ws.Range("I134").Value = 768
$ws.Range("K134").Value = 2304
$ws.Range("M134").Value = 231
# Row 136
$ws.Range("H136").Value = 1431841.8
$ws.Range("I136").Value = 2502699.2
$ws.Range("K136").Value = 7508097.600000001
$ws.Range("M136").Value = -7505547.600000001

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 11155024
$ws.Range("I131").Value = 13890948
$ws.Range("J131").Value = 10354266
$ws.Range("K131").Value = 41672844
$ws.Range("L131").Value = 31062798
$ws.Range("M131").Value = -41667804
$ws.Range("N131").Value = -31072878

$ws = $wb.Worksheets.Item("GSM")
# Row 54
$ws.Range("H54").Value = 12499.5
$ws.Range("J54").Value = 12499.5
$ws.Range("L54").Value = 12499.5
$ws.Range("N54").Value = -13279.5
# Row 70
$ws.Range("H70").Value = 4769304.5
$ws.Range("I70").Value = 6809149.5
$ws.Range("K70").Value = 6809149.5
$ws.Range("M70").Value = -6808879.5
# Row 73
$ws.Range("H73").Value = 4769304.5
$ws.Range("I73").Value = 6809149.5
$ws.Range("K73").Value = 6809149.5
$ws.Range("M73").Value = -6808213.5
# Row 102
$ws.Range("H102").Value = 677695.8
$ws.Range("I102").Value = 1012343.7
$ws.Range("K102").Value = 1012343.7
$ws.Range("M102").Value = -1010721.7
# Row 122
$ws.Range("H122").Value = 692753.7
$ws.Range("I122").Value = 1102406.1
$ws.Range("K122").Value = 3307218.3
$ws.Range("M122").Value = -3304768.3
# Row 126
$ws.Range("H126").Value = 3696.75
$ws.Range("I126").Value = 2061.611
$ws.Range("K126").Value = 6184.833
$ws.Range("M126").Value = -3714.833
# Row 132
$ws.Range("H132").Value = 4242.9443
$ws.Range("I132").Value = 3585.75
$ws.Range("K132").Value = 10757.25
$ws.Range("M132").Value = -8227.25

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6993.125
$ws.Range("I7").Value = 6789
$ws.Range("K7").Value = 6789
$ws.Range("M7").Value = -6677
# Row 40
$ws.Range("H40").Value = 54170668
$ws.Range("I40").Value = 25003250
$ws.Range("J40").Value = 83338090
$ws.Range("K40").Value = 25003250
$ws.Range("L40").Value = 83338090
$ws.Range("M40").Value = -25003114
$ws.Range("N40").Value = -83338362
# Row 41
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
# Row 82
$ws.Range("H82").Value = 4465888.5
$ws.Range("I82").Value = 6250644
$ws.Range("K82").Value = 6250644
$ws.Range("M82").Value = -6250283
# Row 85
$ws.Range("H85").Value = 4465888.5
$ws.Range("I85").Value = 6250644
$ws.Range("K85").Value = 6250644
$ws.Range("M85").Value = -6249396
# Row 122
$ws.Range("H122").Value = 333336000
$ws.Range("I122").Value = 1000000000
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 3000000000
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2999997550
$ws.Range("N122").Value = -16900
# Row 126
$ws.Range("H126").Value = 6993.125
$ws.Range("I126").Value = 6789
$ws.Range("K126").Value = 20367
$ws.Range("M126").Value = -17897
# Row 132
$ws.Range("H132").Value = 5802
$ws.Range("I132").Value = 5110.275
$ws.Range("J132").Value = 6756.1035
$ws.Range("K132").Value = 15330.825
$ws.Range("L132").Value = 20268.3105
$ws.Range("M132").Value = -12800.825
$ws.Range("N132").Value = -25328.3105

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2673.3057
$ws.Range("J122").Value = 2896.3333
$ws.Range("L122").Value = 8688.999899999999
$ws.Range("N122").Value = -13588.9999
# Row 132
$ws.Range("H132").Value = 14622115
$ws.Range("I132").Value = 1589548.2
$ws.Range("J132").Value = 166668740
$ws.Range("K132").Value = 4768644.6
$ws.Range("L132").Value = 500006220
$ws.Range("M132").Value = -4766114.6
$ws.Range("N132").Value = -500011280
# Row 136
$ws.Range("H136").Value = 8925.763000000001
$ws.Range("I136").Value = 2100.3333
$ws.Range("K136").Value = 6300.999899999999
$ws.Range("M136").Value = -3750.999899999999
